$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RESOURCES")

# Insert a new row above the current row 2 (Natural Gas), shifting the
# existing rows (Natural Gas, Electricity, Solar) down by one.
$ws.Rows.Item(2).Insert()

# Give the new row the same formatting as the row right below it (which
# carries the standard data-row style for this table).
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)

# Column F in the data rows uses a bordered style; row 2's F cell should
# stay empty but keep the header-like bordered/filled look used in F1.
$ws.Range("F1").Copy()
$ws.Range("F2").PasteSpecial(-4122)

# Fill in the new "none" resource row.
$ws.Range("A2").Value2 = "none"
$ws.Range("B2").Value2 = "NONE"
$ws.Range("C2").Value2 = 0
$ws.Range("D2").Value2 = 0
$ws.Range("E2").Value2 = 0
$ws.Range("F2").ClearContents()

$excel.CutCopyMode = 0

# Match the selection left behind by the edit.
$ws.Range("A2:F2").Select() | Out-Null
